$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product names
$ws.Range("A1").Value = "Chocolate chips"
$ws.Range("A2").Value = "Crispy bisc"

# Update quantities
$ws.Range("B1").Value = 1300
$ws.Range("B2").Value = 233124

# Update selection to B2
$ws.Range("B2").Select()
